# Append one new log row to each of the four worksheets, following the
# exact same layout (columns A..I) used by the existing rows.
#
# Columns: A=time  B=总长  C=ID  D=实际长度  E=和校验
#          F=总长_DEC  G=ID_DEC  H=实际长度_DEC  I=和校验_DEC

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        [string]$SheetName,
        [int]$RowIndex,
        [double]$TimeValue,
        [string]$TotalLen,
        [string]$Id,
        [string]$ActualLen,
        [string]$Checksum,
        [int]$TotalLenDec,
        $IdDec,
        [int]$ActualLenDec,
        [int]$ChecksumDec,
        [bool]$IdDecIsText
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Cells.Item($RowIndex, 1).Value = $TimeValue
    $ws.Cells.Item($RowIndex, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($RowIndex, 2).Value = $TotalLen
    $ws.Cells.Item($RowIndex, 3).Value = $Id
    $ws.Cells.Item($RowIndex, 4).Value = $ActualLen
    $ws.Cells.Item($RowIndex, 5).Value = $Checksum
    $ws.Cells.Item($RowIndex, 6).Value = $TotalLenDec

    if ($IdDecIsText) {
        # The digit string is too long to round-trip through a double
        # (24 significant digits), so it has to stay text, matching the
        # other rows already on this sheet. Force text entry by switching
        # the cell to the "@" number format before assigning the value,
        # then restore the plain/default style (borrowed from a neighbouring
        # plain-styled cell) so no spurious formatting sticks to the cell.
        $ws.Cells.Item($RowIndex, 7).NumberFormat = "@"
        $ws.Cells.Item($RowIndex, 7).Value = $IdDec
        $ws.Cells.Item($RowIndex, 7).Style = $ws.Cells.Item($RowIndex, 6).Style
    } else {
        $ws.Cells.Item($RowIndex, 7).Value = $IdDec
    }

    $ws.Cells.Item($RowIndex, 8).Value = $ActualLenDec
    $ws.Cells.Item($RowIndex, 9).Value = $ChecksumDec
}

# Note: this interpreter's named-parameter binding (-Param value) does not
# reliably pass values through, so all calls below use positional args.

$idDecNumeric = [double]"5.68631262647114e+23"

# ROW50-FE-LIFTER: new row 43
Add-LogRow "ROW50-FE-LIFTER" 43 45745.67887212963 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x6a" "0xe" 400 $idDecNumeric 362 14 $false

# ROW50-MID-LIFTER: new row 45 (this sheet stores ID_DEC as text, matching
# the existing rows already on the sheet)
Add-LogRow "ROW50-MID-LIFTER" 45 45745.64622685185 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x72" "0x19" 400 "568631262647113771663628" 370 25 $true

# ROW11-FE-LIFTER: new row 43
Add-LogRow "ROW11-FE-LIFTER" 43 45745.69933636574 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x6a" "0x14" 400 $idDecNumeric 362 20 $false

# ROW11-MID-LIFTER: new row 43
Add-LogRow "ROW11-MID-LIFTER" 43 45745.8397540162 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x72" "0x19" 400 $idDecNumeric 370 25 $false
